$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range's last row from column A (row 1 is the header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $d = $dCell.Value2
    $e = $eCell.Value2
    $f = $fCell.Value2

    if ($d -eq $null -or $e -eq $null -or $f -eq $null) {
        continue
    }

    $fStr = [string]([int]$f)
    if ($fStr.Length -ne 8) {
        # Malformed start-date value (e.g. data-entry typo) - leave untouched.
        continue
    }

    $year = [int]$fStr.Substring(0, 4)
    $month = [int]$fStr.Substring(4, 2)
    $day = [int]$fStr.Substring(6, 2)
    $startDate = Get-Date -Year $year -Month $month -Day $day -Hour 0 -Minute 0 -Second 0

    $newE = $e - 1

    if ($newE -le 0) {
        # Stock ran out: restock starting from the previous expiry date,
        # with the remaining-days counter reset to the full duration.
        $expiryDate = $startDate.AddDays([double]$d)
        $fCell.Value = [int]$expiryDate.ToString("yyyyMMdd")
        $eCell.Value = $d
    } else {
        $eCell.Value = $newE
    }
}
